$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 226; existing rows 226:254 shift down to 227:255
$ws.Rows(226).Insert()

# Populate the newly inserted row 226 with its data
$ws.Cells.Item(226, 1).Value = 8
$ws.Cells.Item(226, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(226, 3).Value = "Coquimbo"
$ws.Cells.Item(226, 4).Value = 44491
$ws.Cells.Item(226, 5).Value = 4
$ws.Cells.Item(226, 6).Value = 100114001
$ws.Cells.Item(226, 7).Value = "Papa"
$ws.Cells.Item(226, 8).Value = "Cardinal"
$ws.Cells.Item(226, 9).Value = "1a (cosecha)"
$ws.Cells.Item(226, 10).Value = 3000
$ws.Cells.Item(226, 11).Value = 12000
$ws.Cells.Item(226, 12).Value = 13000
$ws.Cells.Item(226, 13).Value = 12500
$ws.Cells.Item(226, 14).Value = "`$/saco 25 kilos"
$ws.Cells.Item(226, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(226, 16).Value = 500
$ws.Cells.Item(226, 17).Value = 25
$ws.Cells.Item(226, 18).Value = "Hortaliza"
